# Fill in the account-number / amount table on Sheet1 (rows 4-13, cols B:C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (bold) ---
$ws.Range("B4").Value = "Account Number"
$ws.Range("C4").Value = "Amount"
$ws.Range("B4:C4").Font.Bold = $true

# --- Account numbers (column B) and amounts (column C) ---
# First two account numbers were typed with a leading apostrophe (forces
# text + quote-prefix flag) and right aligned.
$ws.Range("B5").Value = "'A0000"
$ws.Range("C5").Value = 785

$ws.Range("B6").Value = "'A1111"
$ws.Range("C6").Value = 6973

$ws.Range("B5:B6").HorizontalAlignment = -4152

# Remaining account numbers are plain text (non-numeric, so already stored
# as text) that are just right aligned - no quote-prefix.
$ws.Range("B7").Value = "A2222"
$ws.Range("C7").Value = 173

$ws.Range("B8").Value = "A3333"
$ws.Range("C8").Value = 0

$ws.Range("B9").Value = "A4444"
$ws.Range("C9").Value = 2984

$ws.Range("B10").Value = "A5555"
$ws.Range("C10").Value = 5059

$ws.Range("B11").Value = "A7777"
$ws.Range("C11").Value = 0

$ws.Range("B12").Value = "A8888"
$ws.Range("C12").Value = 0

$ws.Range("B13").Value = "A9999"
$ws.Range("C13").Value = 2992

$ws.Range("B7:B13").HorizontalAlignment = -4152

# --- Column widths (best-fit sized around the typed data) ---
$ws.Columns("C").ColumnWidth = 15.166666666666666
$ws.Columns("D").ColumnWidth = 7.333333333333333

# --- Leave the cursor where the user stopped typing ---
$ws.Range("F14").Select()
